# BillingSuite.xlsx edit
#
# - TestCases sheet: rename/shuffle the first few test rows
#     A2 "AddStock"     -> "Test1"          (B2 stays "N")
#     A3 "ViewBillTest" -> "Test2", B3 "Y" -> "N"
#     A4 "Test3"        -> "ViewBillTest", B4 "N" -> "Y"
# - Data sheet: the "AddStock" section header becomes "Test2"
# - Selection/active-sheet state moves from TestCases!B2 to Data!A7,
#   with the Data sheet becoming the active tab.

$wb = $excel.ActiveWorkbook
$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsData = $wb.Worksheets.Item("Data")

# --- TestCases sheet: update test-case names / run flags ---
$wsTestCases.Range("A2").Value = "Test1"
$wsTestCases.Range("A3").Value = "Test2"
$wsTestCases.Range("B3").Value = "N"
$wsTestCases.Range("A4").Value = "ViewBillTest"
$wsTestCases.Range("B4").Value = "Y"

# --- Data sheet: rename the AddStock section header to Test2 ---
$wsData.Range("A7").Value = "Test2"

# --- Update selection on TestCases, then switch the active tab/selection to Data ---
$wsTestCases.Range("B11").Select() | Out-Null
$wsData.Activate() | Out-Null
$wsData.Range("A7").Select() | Out-Null
